$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 743.5333000000001
$ws.Range("I33").Value = 123.22222
$ws.Range("J33").Value = 1674
$ws.Range("K33").Value = 123.22222
$ws.Range("L33").Value = 1674
$ws.Range("M33").Value = 105.77778
$ws.Range("N33").Value = -2132
$ws.Range("H129").Value = 893.4
$ws.Range("I129").Value = 646.53845
$ws.Range("J129").Value = 1012.2593
$ws.Range("K129").Value = 1939.61535
$ws.Range("L129").Value = 3036.7779
$ws.Range("M129").Value = 3060.38465
$ws.Range("N129").Value = -13036.7779
$ws.Range("H137").Value = 1303.7457
$ws.Range("I137").Value = 819.7222
$ws.Range("J137").Value = 1516.2439
$ws.Range("K137").Value = 2459.1666
$ws.Range("L137").Value = 4548.7317
$ws.Range("M137").Value = 90.83339999999998
$ws.Range("N137").Value = -9648.7317

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 4333.3335
$ws.Range("J11").Value = 4333.3335
$ws.Range("L11").Value = 4333.3335
$ws.Range("N11").Value = -4621.3335
$ws.Range("H32").Value = 33366.13
$ws.Range("I32").Value = 6814.5386
$ws.Range("J32").Value = 98085.625
$ws.Range("K32").Value = 6814.5386
$ws.Range("L32").Value = 98085.625
$ws.Range("M32").Value = -6527.5386
$ws.Range("N32").Value = -98659.625
$ws.Range("H61").Value = 1724.8718
$ws.Range("I61").Value = 1036.4762
$ws.Range("J61").Value = 2528
$ws.Range("K61").Value = 1036.4762
$ws.Range("L61").Value = 2528
$ws.Range("M61").Value = -824.4762000000001
$ws.Range("N61").Value = -2952
$ws.Range("H74").Value = 1374.9362
$ws.Range("I74").Value = 1038.6957
$ws.Range("J74").Value = 1697.1666
$ws.Range("K74").Value = 1038.6957
$ws.Range("L74").Value = 1697.1666
$ws.Range("M74").Value = -164.6957
$ws.Range("N74").Value = -3445.1666
$ws.Range("H77").Value = 1374.9362
$ws.Range("I77").Value = 1038.6957
$ws.Range("J77").Value = 1697.1666
$ws.Range("K77").Value = 5193.4785
$ws.Range("L77").Value = 8485.833000000001
$ws.Range("M77").Value = -825.4785000000002
$ws.Range("N77").Value = -17221.833
$ws.Range("H136").Value = 1724.8718
$ws.Range("I136").Value = 1036.4762
$ws.Range("J136").Value = 2528
$ws.Range("K136").Value = 3109.4286
$ws.Range("L136").Value = 7584
$ws.Range("M136").Value = -559.4286000000002
$ws.Range("N136").Value = -12684

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 16600
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9826
$ws.Range("H31").Value = 1643.7966
$ws.Range("I31").Value = 1101.8529
$ws.Range("J31").Value = 2380.84
$ws.Range("K31").Value = 1101.8529
$ws.Range("L31").Value = 2380.84
$ws.Range("M31").Value = -806.8529000000001
$ws.Range("N31").Value = -2970.84
$ws.Range("H34").Value = 1643.7966
$ws.Range("I34").Value = 1101.8529
$ws.Range("J34").Value = 2380.84
$ws.Range("K34").Value = 1101.8529
$ws.Range("L34").Value = 2380.84
$ws.Range("M34").Value = -899.8529000000001
$ws.Range("N34").Value = -2784.84
$ws.Range("H64").Value = 40828.332
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H67").Value = 40828.332
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H133").Value = 29000
$ws.Range("J133").Value = 29000
$ws.Range("L133").Value = 29000
$ws.Range("N133").Value = -34060
$ws.Range("H134").Value = 1362.091
$ws.Range("I134").Value = 1297
$ws.Range("J134").Value = 1535.6666
$ws.Range("K134").Value = 3891
$ws.Range("L134").Value = 4606.9998
$ws.Range("M134").Value = -1356
$ws.Range("N134").Value = -9676.9998

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2292.6667
$ws.Range("I70").Value = 470.66666
$ws.Range("J70").Value = 2900
$ws.Range("K70").Value = 1411.99998
$ws.Range("L70").Value = 8700
$ws.Range("M70").Value = -1096.99998
$ws.Range("N70").Value = -9330
$ws.Range("H73").Value = 2292.6667
$ws.Range("I73").Value = 470.66666
$ws.Range("J73").Value = 2900
$ws.Range("K73").Value = 1411.99998
$ws.Range("L73").Value = 8700
$ws.Range("M73").Value = -319.9999800000001
$ws.Range("N73").Value = -10884
$ws.Range("H113").Value = 659.0833
$ws.Range("I113").Value = 655.8333
$ws.Range("J113").Value = 662.3333
$ws.Range("K113").Value = 1967.4999
$ws.Range("L113").Value = 1986.9999
$ws.Range("M113").Value = 202.5001
$ws.Range("N113").Value = -6326.9999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 19612840
$ws.Range("I126").Value = 7496
$ws.Range("J126").Value = 58823530
$ws.Range("K126").Value = 22488
$ws.Range("L126").Value = 176470590
$ws.Range("M126").Value = -20018
$ws.Range("N126").Value = -176475530

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 441.16666
$ws.Range("I22").Value = 410.25
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 410.25
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -115.25
$ws.Range("N22").Value = -1040
$ws.Range("H27").Value = 441.16666
$ws.Range("I27").Value = 410.25
$ws.Range("J27").Value = 450
$ws.Range("K27").Value = 410.25
$ws.Range("L27").Value = 450
$ws.Range("M27").Value = -303.25
$ws.Range("N27").Value = -664
$ws.Range("H132").Value = 6135.7144
$ws.Range("I132").Value = 9200.223
$ws.Range("J132").Value = 3837.3333
$ws.Range("K132").Value = 27600.669
$ws.Range("L132").Value = 11511.9999
$ws.Range("M132").Value = -25070.669
$ws.Range("N132").Value = -16571.9999
$ws.Range("H136").Value = 1079.5758
$ws.Range("I136").Value = 885.96295
$ws.Range("J136").Value = 1950.8334
$ws.Range("K136").Value = 2657.88885
$ws.Range("L136").Value = 5852.5002
$ws.Range("M136").Value = -107.8888499999998
$ws.Range("N136").Value = -10952.5002

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 737.13336
$ws.Range("I113").Value = 367
$ws.Range("K113").Value = 1101
$ws.Range("M113").Value = 1069
$ws.Range("H126").Value = 2316.7273
$ws.Range("I126").Value = 2515.8
$ws.Range("K126").Value = 7547.400000000001
$ws.Range("M126").Value = -5077.400000000001
$ws.Range("H132").Value = 1845.5
$ws.Range("I132").Value = 1465.5
$ws.Range("K132").Value = 4396.5
$ws.Range("M132").Value = -1866.5
